$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.91
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("V4").Value = 1.62
$ws.Range("O10").Value = 1.14
$ws.Range("U10").Value = 2.3
$ws.Range("V10").Value = 1.59
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 2.35
$ws.Range("Q16").Value = 1.67
$ws.Range("R16").Value = 2.15
$ws.Range("Q19").Value = 1.75
$ws.Range("R19").Value = 2.05
$ws.Range("G22").Value = 1.33
$ws.Range("J22").Value = 1.73
$ws.Range("N22").Value = 26
$ws.Range("Z22").Value = 10
$ws.Range("AC22").Value = 26
$ws.Range("AD22").Value = 12
$ws.Range("AH22").Value = 29
$ws.Range("AX22").Value = 34
$ws.Range("BA22").Value = 101
$ws.Range("BB22").Value = 151
$ws.Range("G29").Value = 3.25
$ws.Range("H29").Value = 3.35
$ws.Range("I29").Value = 2.1
$ws.Range("J29").Value = 3.75
$ws.Range("L29").Value = 2.67
$ws.Range("N29").Value = 7.2
$ws.Range("O29").Value = 1.31
$ws.Range("P29").Value = 3.15
$ws.Range("Q29").Value = 1.93
$ws.Range("R29").Value = 1.8
$ws.Range("S29").Value = 1.39
$ws.Range("T29").Value = 2.77
$ws.Range("U29").Value = 1.75
$ws.Range("V29").Value = 1.95
$ws.Range("W29").Value = 9.75
$ws.Range("X29").Value = 17
$ws.Range("Y29").Value = 11.25
$ws.Range("Z29").Value = 45
$ws.Range("AA29").Value = 28
$ws.Range("AB29").Value = 35
$ws.Range("AC29").Value = 7.2
$ws.Range("AD29").Value = 6.4
$ws.Range("AE29").Value = 14
$ws.Range("AF29").Value = 65
$ws.Range("AG29").Value = 500
$ws.Range("AH29").Value = 7.5
$ws.Range("AI29").Value = 10
$ws.Range("AJ29").Value = 8.75
$ws.Range("AK29").Value = 19.5
$ws.Range("AL29").Value = 17
$ws.Range("AM29").Value = 28
$ws.Range("AN29").Value = 5.2
$ws.Range("AO29").Value = 17.5
$ws.Range("AP29").Value = 24
$ws.Range("AQ29").Value = 90
$ws.Range("AR29").Value = 120
$ws.Range("AS29").Value = 300
$ws.Range("AT29").Value = 2.77
$ws.Range("AU29").Value = 7
$ws.Range("AV29").Value = 60
$ws.Range("AW29").Value = 4.05
$ws.Range("AX29").Value = 10.75
$ws.Range("AZ29").Value = 40
$ws.Range("BA29").Value = 70
$ws.Range("M33").Value = 1.07
$ws.Range("N33").Value = 9
$ws.Range("Q33").Value = 2.2
$ws.Range("R33").Value = 1.65
